$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 15:35"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1647043
$ws.Range("C4").Value = 1949
$ws.Range("E4").Value = 1146044
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 97687

# Alemania (row 11)
$ws.Range("B11").Value = 179758
$ws.Range("C11").Value = 45
$ws.Range("E11").Value = 11505
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 8353

# India (row 14)
$ws.Range("B14").Value = 127358
$ws.Range("C14").Value = 2564
$ws.Range("E14").Value = 71327
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = 3759

# Arabia Saudita (row 18)
$ws.Range("B18").Value = 70161
$ws.Range("C18").Value = 2442
$ws.Range("D18").Value = 41236
$ws.Range("E18").Value = 28546
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 379

# Catar (row 24)
$ws.Range("B24").Value = 42213
$ws.Range("C24").Value = 1732
$ws.Range("D24").Value = 8513
$ws.Range("E24").Value = 33679
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 21

# Suiza (row 30)
$ws.Range("E30").Value = 920
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 1905

# Chequia (row 53)
$ws.Range("B53").Value = 8820
$ws.Range("C53").Value = 7
$ws.Range("D53").Value = 6026
$ws.Range("E53").Value = 2482

# Noruega (row 55)
$ws.Range("B55").Value = 8340
$ws.Range("C55").Value = 8
$ws.Range("E55").Value = 378

# Mozambique (row 158)
$ws.Range("B158").Value = 168
$ws.Range("C158").Value = 4
$ws.Range("E158").Value = 120
